$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: months/values, row 13 added, labels & values shifted forward (2015/2016 -> 2018/2019)
$labels = @("October 2018", "November 2018", "December 2018", "January 2019", "February 2019", "March 2019", "April 2019", "May 2019", "June 2019", "July 2019", "August 2019", "September 2019", "October 2019")
$values = @(497, 373, 568, 554, 581, 670, 612, 629, 371, 691, 404, 563, 519)

# Force column A to be treated as text so month names aren't auto-converted to dates
$ws.Range("A1:A13").NumberFormat = "@"

for ($i = 0; $i -lt $labels.Count; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $labels[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Update chart style and extend the series range to include the new row 13
$co = $ws.ChartObjects().Item(1)
$chart = $co.Chart
$chart.ChartStyle = 2

$ser = $chart.SeriesCollection(1)
$ser.Formula = "=SERIES(,Sheet1!`$A`$1:`$A`$13,Sheet1!`$B`$1:`$B`$13,1)"
